# Adding a more simple tab for testing purposes
#
# Create a new worksheet "Book_02" positioned between "Book_01" and
# "Dummy_Tab". It is built as a trimmed-down copy of "Book_01": same
# header row, same first data row, and the last data row of Book_01
# kept as its final row, with its own (smaller) AutoFilter range.

$wb = $excel.ActiveWorkbook
$book1 = $wb.Worksheets.Item("Book_01")
$dummyTab = $wb.Worksheets.Item("Dummy_Tab")

# Copy Book_01 so the new sheet inherits identical page/view/column
# formatting, then drop it right before Dummy_Tab and rename it.
$book1.Copy($dummyTab)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "Book_02"

# Keep header (row 1) + first data row (row 2), drop everything else
# except the very last data row, which we will restore afterwards.
$newSheet.Range("A3:A59").EntireRow.Delete()

# Re-apply AutoFilter while the sheet only has 2 rows so the filter
# range stays B1:B2 instead of snapping to a larger block.
$newSheet.AutoFilterMode = $false
$newSheet.Range("B1:B2").AutoFilter()

# Bring back Book_01's last row (A058 / Zimbabyoue / ...) as row 3,
# copying both values and formatting from Book_01's row 2 style.
$book1.Range("A2:E2").Copy()
$newSheet.Range("A3:E3").PasteSpecial(-4122)
$newSheet.Range("A3").Value = "A058"
$newSheet.Range("B3").Value = "Zimbabyoue"
$newSheet.Range("C3").Value = "9212 Lokeren, Belgium"
$newSheet.Range("D3").Value = "erg"
$newSheet.Range("E3").Value = "erg"

# Register the hidden _FilterDatabase name for the new sheet, same as
# Book_01 already has for itself.
$fd = $newSheet.Names.Add("_xlnm._FilterDatabase", "=Book_02!`$B`$1:`$B`$2")
foreach ($nm in $wb.Names) {
    if ($nm.Name -like "*_FilterDatabase*" -and $nm.RefersTo -like "*Book_02*") {
        $nm.Visible = $false
    }
}

# Match the saved selection on the new tab.
$newSheet.Range("B10").Select()
